# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# for the cryptos worksheet, matching the provided OOXML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.858.21"
$ws.Range("E2").Value = "  +5.88%  "
$ws.Range("D3").Value = "3.069.12"
$ws.Range("E3").Value = "  +3.41%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.75"
$ws.Range("E5").Value = "  +3.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.96"
$ws.Range("E6").Value = "  +5.98%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "3.051.62"
$ws.Range("E8").Value = "  +2.94%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.525"
$ws.Range("E9").Value = "  +2.24%  "
$ws.Range("E10").Value = "  +6.46%  "
$ws.Range("E11").Value = "  +14.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.464"
$ws.Range("E12").Value = "  +2.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000238"
$ws.Range("E13").Value = "  +6.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.76"
$ws.Range("E14").Value = "  +4.37%  "
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("D16").Value = "3.573.80"
$ws.Range("E16").Value = "  +3.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.22"
$ws.Range("E17").Value = "  +4.34%  "
$ws.Range("D18").Value = "3.057.76"
$ws.Range("E18").Value = "  +2.80%  "
$ws.Range("D19").Value = "61.795.31"
$ws.Range("E19").Value = "  +5.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "449.41"
$ws.Range("E20").Value = "  +6.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.88"
$ws.Range("E21").Value = "  +3.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.732"
$ws.Range("E22").Value = "  +3.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.28"
$ws.Range("E23").Value = "  +3.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.68"
$ws.Range("E24").Value = "  +2.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.90"
$ws.Range("E25").Value = "  +2.61%  "
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("E27").Value = "  +7.82%  "
$ws.Range("E28").Value = "  -0.23%  "
$ws.Range("E29").Value = "  +5.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.10"
$ws.Range("E30").Value = "  +5.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.47"
$ws.Range("E31").Value = "  +7.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.52"
$ws.Range("E32").Value = "  +3.84%  "
$ws.Range("E33").Value = "  +8.37%  "
$ws.Range("D34").Value = "0.0₃0811"
$ws.Range("E34").Value = "  +8.74%  "
$ws.Range("E35").Value = "  +4.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.06"
$ws.Range("E36").Value = "  +6.36%  "
$ws.Range("E37").Value = "  +6.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "50.00"
$ws.Range("E38").Value = "  +3.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.97"
$ws.Range("E39").Value = "  +8.94%  "
$ws.Range("E40").Value = "  +1.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "414.19"
$ws.Range("E41").Value = "  +4.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0365"
$ws.Range("E42").Value = "  +5.83%  "
$ws.Range("D43").Value = "2.786.10"
$ws.Range("E43").Value = "  +2.12%  "
$ws.Range("E44").Value = "  +1.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.264"
$ws.Range("E45").Value = "  +9.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "37.82"
$ws.Range("E46").Value = "  +19.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.11"
$ws.Range("E47").Value = "  +6.20%  "
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.87"
$ws.Range("E49").Value = "  -0.75%  "
$ws.Range("E50").Value = "  +2.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.18"
$ws.Range("E51").Value = "  +4.82%  "
